$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# (e.g. "0.9982", "288.19") are preserved exactly as text, matching
# the original inline-string formatting instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '22.168.89'
$ws.Range("E2").Value = '  -1.20%  '
$ws.Range("D3").Value = '1.559.74'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("D4").Value = '0.9982'
$ws.Range("E4").Value = '  -0.38%  '
$ws.Range("D5").Value = '0.9982'
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("D6").Value = '288.19'
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = '0.3798'
$ws.Range("E7").Value = '  +2.46%  '
$ws.Range("D8").Value = '0.3296'
$ws.Range("E8").Value = '  -0.33%  '
$ws.Range("D9").Value = '43.54'
$ws.Range("E9").Value = '  -9.92%  '
$ws.Range("D10").Value = '1.149'
$ws.Range("E10").Value = '  +1.20%  '
$ws.Range("D11").Value = '0.07388'
$ws.Range("E11").Value = '  -1.45%  '
$ws.Range("D12").Value = '0.9982'
$ws.Range("E12").Value = '  -0.40%  '
$ws.Range("D13").Value = '20.27'
$ws.Range("E13").Value = '  -2.13%  '
$ws.Range("D14").Value = '5.840'
$ws.Range("E14").Value = '  -1.51%  '
$ws.Range("D15").Value = '6.862'
$ws.Range("E15").Value = '  -0.16%  '
$ws.Range("D16").Value = '1.554.60'
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("D17").Value = '0.00001112'
$ws.Range("E17").Value = '  -0.57%  '
$ws.Range("D18").Value = '0.06628'
$ws.Range("E18").Value = '  -1.80%  '
$ws.Range("D19").Value = '86.17'
$ws.Range("E19").Value = '  -1.56%  '
$ws.Range("D20").Value = '6.411'
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("D21").Value = '0.9983'
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("D22").Value = '16.13'
$ws.Range("E22").Value = '  -2.95%  '
$ws.Range("D23").Value = '11.74'
$ws.Range("E23").Value = '  -2.45%  '
$ws.Range("D24").Value = '22.167.62'
$ws.Range("E24").Value = '  -1.20%  '
$ws.Range("D25").Value = '2.314'
$ws.Range("E25").Value = '  -3.14%  '
$ws.Range("D26").Value = '2.540'
$ws.Range("E26").Value = '  -1.72%  '
$ws.Range("D27").Value = '150.43'
$ws.Range("E27").Value = '  -2.31%  '
$ws.Range("D28").Value = '19.17'
$ws.Range("E28").Value = '  -2.76%  '
$ws.Range("D29").Value = '4.922'
$ws.Range("E29").Value = '  -1.99%  '
$ws.Range("D30").Value = '121.91'
$ws.Range("E30").Value = '  -2.05%  '
$ws.Range("D31").Value = '1.727.72'
$ws.Range("E31").Value = '  -1.05%  '
$ws.Range("D32").Value = '1.079'
$ws.Range("E32").Value = '  +2.18%  '
$ws.Range("D33").Value = '5.977'
$ws.Range("E33").Value = '  -2.16%  '
$ws.Range("D34").Value = '1.850'
$ws.Range("E34").Value = '  -8.22%  '
$ws.Range("D35").Value = '0.08245'
$ws.Range("E35").Value = '  -1.53%  '
$ws.Range("D36").Value = '9.359'
$ws.Range("E36").Value = '  -4.38%  '
$ws.Range("D37").Value = '0.02343'
$ws.Range("E37").Value = '  -5.10%  '
$ws.Range("D38").Value = '0.06273'
$ws.Range("E38").Value = '  -1.99%  '
$ws.Range("D39").Value = '5.310'
$ws.Range("E39").Value = '  -0.60%  '
$ws.Range("D40").Value = '0.2171'
$ws.Range("E40").Value = '  -3.97%  '
$ws.Range("D41").Value = '1.255'
$ws.Range("E41").Value = '  -2.57%  '
$ws.Range("D42").Value = '11.09'
$ws.Range("E42").Value = '  -1.97%  '
$ws.Range("D43").Value = '0.6093'
$ws.Range("E43").Value = '  -3.57%  '
$ws.Range("D44").Value = '0.9987'
$ws.Range("E44").Value = '  -0.23%  '
$ws.Range("D45").Value = '13.74'
$ws.Range("E45").Value = '  -0.78%  '
$ws.Range("D46").Value = '3.742'
$ws.Range("E46").Value = '  -0.75%  '
$ws.Range("D47").Value = '0.5887'
$ws.Range("E47").Value = '  -4.34%  '
$ws.Range("D48").Value = '2.000'
$ws.Range("E48").Value = '  -2.91%  '
$ws.Range("D49").Value = '122.38'
$ws.Range("E49").Value = '  -2.97%  '
$ws.Range("D50").Value = '1.179'
$ws.Range("E50").Value = '  -2.90%  '
$ws.Range("D51").Value = '0.07036'
$ws.Range("E51").Value = '  -2.48%  '
